# "add some data to dataset" - append one new row (A17:B17) describing a
# new unit test, matching the two new shared strings added to the sheet:
#   A17 = test function name
#   B17 = "what it tests" description (long text -> wraps)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "test_tnorm_matrix_consistency_with_scalar_application"
$ws.Range("B17").Value = "output correctness when matrix (2d np.array) is passed to be sure lb and ub in itfrs are correct"

# The new row's first cell wraps its text (distinct style from the plain
# column-A cells above it).
$ws.Range("A17").WrapText = $true

# Column B (long "what it tests" text) and column D (long "how it works"
# text) grow to better fit the new/expanded content.
$ws.Columns.Item(2).ColumnWidth = 78.83333333333333
$ws.Columns.Item(4).ColumnWidth = 93.5

# Leave the view scrolled/selected where editing ended up.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B21").Select()
